$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.199.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.981.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  -4.34%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.110.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.979.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.96%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.985"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "44.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0357"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.731.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.18%  "
